# Update "想去人数" (want-to-go count) figures in both the "展览" and
# "全部类型" worksheets, which hold identical data tables.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1089
    5  = 3082
    7  = 2489
    9  = 123
    11 = 1262
    15 = 1122
    16 = 312
    18 = 23
    19 = 24
    20 = 104
    21 = 66
    23 = 206
    24 = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
